$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.895.24'
$ws.Range("D3").Value = '1.872.60'
$ws.Range("E3").Value = '  -0.86%  '
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = "'0.7403"
$ws.Range("E5").Value = '  -2.99%  '
$ws.Range("D6").Value = "'242.36"
$ws.Range("E6").Value = '  -0.14%  '
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = "'0.3150"
$ws.Range("E8").Value = '  +0.86%  '
$ws.Range("E9").Value = '  +0.70%  '
$ws.Range("E10").Value = '  -3.77%  '
$ws.Range("D11").Value = "'0.08310"
$ws.Range("E11").Value = '  -2.35%  '
$ws.Range("D12").Value = "'0.7492"
$ws.Range("E12").Value = '  -1.73%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'5.377"
$ws.Range("E13").Value = '  +0.26%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.853.35'
$ws.Range("E14").Value = '  +3.96%  '
$ws.Range("E15").Value = '  -1.57%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '29.918.34'
$ws.Range("E16").Value = '  +0.57%  '
$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D17").Value = "'6.105"
$ws.Range("E17").Value = '  -0.56%  '
$ws.Range("D18").Value = "'246.86"
$ws.Range("E18").Value = '  +1.35%  '
$ws.Range("E19").Value = '  -1.61%  '
$ws.Range("D20").Value = "'0.000007834"
$ws.Range("E20").Value = '  +0.50%  '
$ws.Range("D21").Value = "'0.9997"
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").Value = '2.133.89'
$ws.Range("E22").Value = '  +3.08%  '
$ws.Range("D23").Value = "'7.992"
$ws.Range("E23").Value = '  +0.28%  '
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").Value = "'0.1540"
$ws.Range("E25").Value = '  -4.57%  '
$ws.Range("E26").Value = '  -1.33%  '
$ws.Range("D27").Value = "'165.53"
$ws.Range("E27").Value = '  +2.25%  '
$ws.Range("E28").Value = '  -0.57%  '
$ws.Range("D29").Value = "'2.018"
$ws.Range("E29").Value = '  -0.75%  '
$ws.Range("D30").Value = "'1.493"
$ws.Range("E30").Value = '  +1.44%  '
$ws.Range("D31").Value = "'4.572"
$ws.Range("E31").Value = '  +1.85%  '
$ws.Range("D32").Value = "'1.533"
$ws.Range("E32").Value = '  -0.07%  '
$ws.Range("D33").Value = "'4.216"
$ws.Range("E33").Value = '  +2.98%  '
$ws.Range("D34").Value = "'0.05315"
$ws.Range("E34").Value = '  -2.33%  '
$ws.Range("D35").Value = "'1.237"
$ws.Range("E35").Value = '  -0.34%  '
$ws.Range("D36").Value = "'0.7489"
$ws.Range("E36").Value = '  +0.90%  '
$ws.Range("D37").Value = "'1.001"
$ws.Range("E37").Value = '  +0.23%  '
$ws.Range("D38").Value = "'2.695"
$ws.Range("E38").Value = '  -0.14%  '
$ws.Range("D39").Value = "'0.01960"
$ws.Range("E39").Value = '  +0.77%  '
$ws.Range("E40").Value = '  -1.04%  '
$ws.Range("D41").Value = "'0.4514"
$ws.Range("E41").Value = '  +1.22%  '
$ws.Range("D42").Value = '1.113.01'
$ws.Range("E42").Value = '  +0.71%  '
$ws.Range("D43").Value = "'6.132"
$ws.Range("E43").Value = '  +1.19%  '
$ws.Range("D44").Value = "'72.24"
$ws.Range("E44").Value = '  -0.57%  '
$ws.Range("D45").Value = "'0.8637"
$ws.Range("E45").Value = '  +1.27%  '
$ws.Range("E46").Value = '  +1.35%  '
$ws.Range("E47").Value = '  -0.07%  '
$ws.Range("D48").Value = "'1.859"
$ws.Range("E48").Value = '  -0.31%  '
$ws.Range("D49").Value = "'7.615"
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("D50").Value = "'9.514"
$ws.Range("E50").Value = '  -2.17%  '
$ws.Range("D51").Value = '2.032.75'
$ws.Range("E51").Value = '  +1.67%  '
